$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new data rows. Insert two blank rows right before the
# existing row 35 so the old rows 35-38 shift down to become rows 37-40
# (their contents are unchanged by this edit), then populate the two new
# rows (35 and 36) with the "new" records, and tweak the Fecha/Volumen
# values on rows 33 and 34.
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()

# Row 33: Fecha 44985 -> 44988, Volumen 50 -> 30 (everything else unchanged)
$ws.Range("D33").Value = 44988
$ws.Range("M33").Value = 30

# Row 34: Fecha 44985 -> 44988, Volumen 50 -> 30 (everything else unchanged)
$ws.Range("D34").Value = 44988
$ws.Range("M34").Value = 30

# New row 35
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 44985
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100101
$ws.Range("H35").Value = "Berries"
$ws.Range("I35").Value = 100101001
$ws.Range("J35").Value = "Arándano (blue)"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 50
$ws.Range("N35").Value = 3000
$ws.Range("O35").Value = 3000
$ws.Range("P35").Value = 3000
$ws.Range("Q35").Value = "$/bandeja 2 kilos"
$ws.Range("R35").Value = "Provincia de Diguillín"
$ws.Range("S35").Value = 1500
$ws.Range("T35").Value = 2

# New row 36
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44985
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100101
$ws.Range("H36").Value = "Berries"
$ws.Range("I36").Value = 100101001
$ws.Range("J36").Value = "Arándano (blue)"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 50
$ws.Range("N36").Value = 2500
$ws.Range("O36").Value = 2500
$ws.Range("P36").Value = 2500
$ws.Range("Q36").Value = "$/bandeja 2 kilos"
$ws.Range("R36").Value = "Provincia de Diguillín"
$ws.Range("S36").Value = 1250
$ws.Range("T36").Value = 2

Write-Host "Edit applied"
